# Update the "取得日時" (acquisition timestamp) column on the "ランサーズ" sheet
# for all existing data rows (rows 2-18) from "2026-02-06 18:45:41" to
# "2026-02-06 18:56:52" — i.e. a new scrape/append pass was run at 18:56 JST.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2026-02-06 18:45:41"
$newTimestamp = "2026-02-06 18:56:52"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 18
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value2 = $newTimestamp
    }
}
